$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 753, shifting existing rows (753-794) down to (754-795)
$ws.Range("A753").EntireRow.Insert()

# Populate the newly inserted row 753 with the new data point.
# Date-like text ("2026/01/31") would otherwise be auto-converted by Excel
# into a real date serial number, so prefix with an apostrophe to force
# plain text entry (same trick used in the Excel UI), matching the
# original sheet's inlineStr/shared-string date-as-text convention.
$ws.Range("A753").Value = "'2026/01/31"
$ws.Range("A753").Style = "Normal"
$ws.Range("B753").Value = "土"
$ws.Range("C753").Value = 3
$ws.Range("D753").Value = 201
